# Skills data import: add a "type" column and two new skills
# (Alchemy, Disenchanting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "type" column as column B, shifting the rest right ---
$ws.Columns("B").Insert()

# Header
$ws.Range("B1").Value = "type"

# Type values for the existing 11 skills (rows 2-12)
$types = @(0, 0, 0, 1, 1, 1, 1, 1, 2, 5, 6)
for ($i = 0; $i -lt $types.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $types[$i]
}

# --- Add the two new skills: Alchemy (row 13) and Disenchanting (row 14) ---
$ws.Range("A13").Value = "Alchemy"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "Alchemy is a skill that is used by you crafting new and magical items that can do one of two things: Devastating damage to kingdoms, or give you boons that can stack and last for a set number of hours."
$ws.Range("D13").Value = 400
$ws.Range("L13").Value = 0.0025

$ws.Range("A14").Value = "Disenchanting"
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = "This is used by you disenchanting items that either popup in chat as drops or in your inventory by clicking Disenchant All.`nShould you fail to disenchant an item, you'll only get 1 Gold Dust. Should you succeed you can get between 1 and 150 Gold Dust, where as destroying only gets you 1-25 gold dust and no disenchanting experience."
$ws.Range("D14").Value = 999
$ws.Range("L14").Value = 0.001

# Wrap the (multi-line) Disenchanting description
$ws.Range("C14").WrapText = $true

# --- Column widths: shift former widths right by one, give the new
#     "type" column (B) a width of 5, and the new trailing column (M) 25 ---
# (Excel's ColumnWidth setter pads by 5/6 of a character vs. the raw OOXML
# width, so we compensate to land on the exact target widths.)
$ws.Columns("A").ColumnWidth = 23 - 5/6
$ws.Columns("B").ColumnWidth = 5 - 5/6
$ws.Columns("C").ColumnWidth = 930 - 5/6
$ws.Columns("D").ColumnWidth = 11 - 5/6
$ws.Columns("E").ColumnWidth = 37 - 5/6
$ws.Columns("F").ColumnWidth = 38 - 5/6
$ws.Columns("G").ColumnWidth = 32 - 5/6
$ws.Columns("H").ColumnWidth = 41 - 5/6
$ws.Columns("I").ColumnWidth = 39 - 5/6
$ws.Columns("J").ColumnWidth = 28 - 5/6
$ws.Columns("K").ColumnWidth = 11 - 5/6
$ws.Columns("L").ColumnWidth = 25 - 5/6
$ws.Columns("M").ColumnWidth = 25 - 5/6

# Re-fit row 14 now that column C is back to its (huge) final width, so the
# wrapped two-line description doesn't leave a stale/custom row height.
$ws.Rows(14).AutoFit()

# --- Selection ends on the newly added description cell ---
$ws.Range("C14").Select()
